# "aggiornamento fino a 28 luglio" - append new daily rows (302-328) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(44376, 0, 3, 18.78522229179712),
    @(44377, 0, 3, 18.78522229179712),
    @(44378, 0, 3, 18.78522229179712),
    @(44379, 0, 3, 18.78522229179712),
    @(44380, 0, 1, 6.261740763932373),
    @(44381, 0, 1, 6.261740763932373),
    @(44382, 0, 0, 0),
    @(44383, 0, 0, 0),
    @(44384, 0, 0, 0),
    @(44385, 0, 0, 0),
    @(44386, 0, 0, 0),
    @(44387, 0, 0, 0),
    @(44388, 0, 0, 0),
    @(44389, 0, 0, 0),
    @(44390, 0, 0, 0),
    @(44391, 1, 1, 6.261740763932373),
    @(44392, 1, 2, 12.52348152786475),
    @(44393, 0, 2, 12.52348152786475),
    @(44394, 1, 3, 18.78522229179712),
    @(44395, 0, 3, 18.78522229179712),
    @(44396, 0, 3, 18.78522229179712),
    @(44397, 0, 3, 18.78522229179712),
    @(44398, 0, 2, 12.52348152786475),
    @(44399, 4, 5, 31.30870381966186),
    @(44400, 0, 5, 31.30870381966186),
    @(44401, 1, 5, 31.30870381966186),
    @(44402, 0, 5, 31.30870381966186)
)

$startRow = 302
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Column A in the existing data carries the date-centered/bold/bordered style (same
# as the header row of dates); replicate it onto the newly appended date cells.
$ws.Range("A301").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
